# Updates cryptos list values (price/volume, and a row reorder for LEO / InternetComputer)
# to match the latest scrape, per commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "66.443.97"
$ws.Range("E2").Value = "  +3.75%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.503.10"
$ws.Range("E3").Value = "  +4.40%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  -0.08%  "

# Row 5: BNB
$ws.Range("D5").Value = "'558.15"
$ws.Range("E5").Value = "  +6.69%  "

# Row 6: Solana
$ws.Range("D6").Value = "'185.29"
$ws.Range("E6").Value = "  +6.79%  "

# Row 7: XRP
$ws.Range("D7").Value = "'0.636"
$ws.Range("E7").Value = "  +7.24%  "

# Row 8: LidoStakedEther
$ws.Range("D8").Value = "3.496.38"
$ws.Range("E8").Value = "  +4.27%  "

# Row 9: USDC
$ws.Range("E9").Value = "  -0.03%  "

# Row 10: Cardano
$ws.Range("D10").Value = "'0.632"
$ws.Range("E10").Value = "  +4.17%  "

# Row 11: Dogecoin
$ws.Range("D11").Value = "'0.153"
$ws.Range("E11").Value = "  +14.31%  "

# Row 12: Avalanche
$ws.Range("D12").Value = "'54.48"
$ws.Range("E12").Value = "  +2.32%  "

# Row 13: ShibaInu
$ws.Range("E13").Value = "  +5.81%  "

# Row 14: Polkadot
$ws.Range("E14").Value = "  +2.72%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.062.81"
$ws.Range("E15").Value = "  +3.67%  "

# Row 16: WrappedEther
$ws.Range("D16").Value = "3.499.67"
$ws.Range("E16").Value = "  +3.86%  "

# Row 17: Chainlink
$ws.Range("D17").Value = "'18.66"
$ws.Range("E17").Value = "  +6.48%  "

# Row 18: TRON
$ws.Range("E18").Value = "  +3.29%  "

# Row 19: WrappedBTC
$ws.Range("D19").Value = "66.435.98"
$ws.Range("E19").Value = "  +3.89%  "

# Row 20: Uniswap
$ws.Range("E20").Value = "  +7.24%  "

# Row 21: Polygon
$ws.Range("E21").Value = "  +3.41%  "

# Row 22: BitcoinCash
$ws.Range("D22").Value = "'421.90"
$ws.Range("E22").Value = "  +12.76%  "

# Row 23: PancakeSwap
$ws.Range("E23").Value = "  +10.70%  "

# Row 24: Litecoin
$ws.Range("D24").Value = "'86.13"
$ws.Range("E24").Value = "  +5.62%  "

# Row 25: Toncoin
$ws.Range("E25").Value = "  -2.37%  "

# Row 26: RenderToken
$ws.Range("D26").Value = "'10.99"
$ws.Range("E26").Value = "  -4.85%  "

# Row 27: ImmutableX
$ws.Range("E27").Value = "  +7.48%  "

# Row 28: LEO
$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "'12.25"
$ws.Range("E28").Value = "  +8.60%  "

# Row 29: InternetComputer(DFINITY)
$ws.Range("B29").Value = "LEO"
$ws.Range("C29").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D29").Value = "'6.08"
$ws.Range("E29").Value = "  -1.45%  "

# Row 30: Filecoin
$ws.Range("D30").Value = "'9.11"
$ws.Range("E30").Value = "  +11.03%  "

# Row 31: EthereumClassic
$ws.Range("D31").Value = "'30.20"
$ws.Range("E31").Value = "  +4.73%  "

# Row 32: Bittensor
$ws.Range("D32").Value = "'626.84"
$ws.Range("E32").Value = "  -0.21%  "

# Row 33: NEARProtocol
$ws.Range("D33").Value = "'6.59"
$ws.Range("E33").Value = "  +2.46%  "

# Row 34: Cosmos
$ws.Range("D34").Value = "'11.72"
$ws.Range("E34").Value = "  +4.79%  "

# Row 35: Hedera
$ws.Range("E35").Value = "  +4.99%  "

# Row 36: OKB
$ws.Range("D36").Value = "'60.15"
$ws.Range("E36").Value = "  +3.93%  "

# Row 37: PEPE
$ws.Range("D37").Value = "0.0₃0824"
$ws.Range("E37").Value = "  +12.59%  "

# Row 38: Kaspa
$ws.Range("D38").Value = "'0.147"
$ws.Range("E38").Value = "  +18.18%  "

# Row 39: InjectiveProtocol
$ws.Range("D39").Value = "'37.88"
$ws.Range("E39").Value = "  +4.91%  "

# Row 40: Dai
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.22%  "

# Row 41: TheGraph
$ws.Range("D41").Value = "'0.386"
$ws.Range("E41").Value = "  +1.83%  "

# Row 42: Stacks
$ws.Range("D42").Value = "'3.42"
$ws.Range("E42").Value = "  +12.96%  "

# Row 43: Maker
$ws.Range("D43").Value = "3.119.39"
$ws.Range("E43").Value = "  +5.03%  "

# Row 44: FirstDigitalUSD
$ws.Range("E44").Value = "  -0.11%  "

# Row 45: Fetch.AI
$ws.Range("D45").Value = "'2.61"
$ws.Range("E45").Value = "  -2.20%  "

# Row 46: ThetaToken
$ws.Range("E46").Value = "  +9.70%  "

# Row 47: ApeXProtocol
$ws.Range("D47").Value = "'3.35"
$ws.Range("E47").Value = "  +11.02%  "

# Row 48: VeChain
$ws.Range("E48").Value = "  +5.10%  "

# Row 49: WEMIXToken
$ws.Range("E49").Value = "  +2.82%  "

# Row 50: Stellar
$ws.Range("D50").Value = "'0.134"
$ws.Range("E50").Value = "  +7.09%  "

# Row 51: Monero
$ws.Range("D51").Value = "'139.15"
$ws.Range("E51").Value = "  +2.40%  "
